# added ifo gdp component analysis preprocessing
# Updates the ifoCAST error series (GVA) data sheet with refreshed values
# for the latest vintages, including two newly added data points
# (D19 and C20) that extend the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 0.07541333303608777

$ws.Range("C6").Value = 1.696990870036088
$ws.Range("D6").Value = 0.44417418258209

$ws.Range("B7").Value = 0.1903804690360879
$ws.Range("C7").Value = 1.06983264558209
$ws.Range("D7").Value = -0.9590976180317647

$ws.Range("B8").Value = 0.66923223158209
$ws.Range("C8").Value = -0.8452485510317647

$ws.Range("B9").Value = -1.085847230031765
$ws.Range("D9").Value = 0.9920349685636995

$ws.Range("C10").Value = 0.9333254135636995
$ws.Range("D10").Value = 0.8577491755291454

$ws.Range("B11").Value = 1.1210615525637
$ws.Range("C11").Value = 0.3885492905291454
$ws.Range("D11").Value = -0.8873425837224346

$ws.Range("B12").Value = 0.4182579295291455
$ws.Range("C12").Value = -0.6454920347224345
$ws.Range("D12").Value = -1.024285167559779

$ws.Range("B13").Value = -0.5544081717224345
$ws.Range("C13").Value = -1.191515643655161
$ws.Range("D13").Value = 1.016949629488724

$ws.Range("B14").Value = -1.26539928353432
$ws.Range("C14").Value = 0.5246197534887239
$ws.Range("D14").Value = 0.3082097950934801

$ws.Range("B15").Value = 0.3295457764887239
$ws.Range("C15").Value = 0.5169714940934801
$ws.Range("D15").Value = 0.8421929667131937

$ws.Range("B16").Value = -0.02970344090651991
$ws.Range("C16").Value = -0.1799067152868063

$ws.Range("B17").Value = -0.2670276532868063
$ws.Range("D17").Value = -0.4617206544916527

$ws.Range("C18").Value = -0.1834133014916527
$ws.Range("D18").Value = 0.5298148361758083

$ws.Range("B19").Value = 0.03353349250834725
$ws.Range("C19").Value = 0.1236970551758083
$ws.Range("D19").Value = 0.5311946523564907

$ws.Range("B20").Value = -0.1256759188241917
$ws.Range("C20").Value = -0.02418658464350926
